# Add two new columns, I ("I0") and J ("IF"), to Sheet1 with header row and
# per-row values for rows 2-66. Matches the diff: dimension grows from
# A1:H66 to A1:J66, headers "I0"/"IF" take on the same style as the other
# header cells (copied from H1), and data rows get literal numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: copy the existing header cell format (bold, bordered, centered)
# from H1 onto I1:J1, then set the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row data: row, I-value, J-value
$data = @(
    @(2, 1, 1),
    @(3, 5, 5),
    @(4, 7, 7),
    @(5, 1, 1),
    @(6, 8, 8),
    @(7, 7, 7),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 5, 5),
    @(11, 9, 9),
    @(12, 5, 5),
    @(13, 5, 6),
    @(14, 8, 8),
    @(15, 3, 3),
    @(16, 7, 7),
    @(17, 7, 7),
    @(18, 6, 6),
    @(19, 7, 7),
    @(20, 6, 7),
    @(21, 6, 6),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 5, 5),
    @(25, 8, 8),
    @(26, 6, 6),
    @(27, 4, 5),
    @(28, 5, 5),
    @(29, 5, 6),
    @(30, 8, 8),
    @(31, 4, 5),
    @(32, 7, 8),
    @(33, 6, 6),
    @(34, 9, 9),
    @(35, 7, 8),
    @(36, 9, 9),
    @(37, 6, 6),
    @(38, 8, 8),
    @(39, 8, 8),
    @(40, 6, 6),
    @(41, 8, 8),
    @(42, 5, 5),
    @(43, 9, 9),
    @(44, 7, 7),
    @(45, 7, 7),
    @(46, 9, 9),
    @(47, 5, 5),
    @(48, 8, 8),
    @(49, 9, 9),
    @(50, 7, 7),
    @(51, 6, 6),
    @(52, 6, 7),
    @(53, 7, 8),
    @(54, 5, 6),
    @(55, 6, 6),
    @(56, 7, 7),
    @(57, 2, 3),
    @(58, 6, 6),
    @(59, 8, 8),
    @(60, 7, 8),
    @(61, 7, 7),
    @(62, 6, 6),
    @(63, 4, 4),
    @(64, 5, 6),
    @(65, 7, 7),
    @(66, 8, 8)
)

foreach ($item in $data) {
    $r = $item[0]
    $iVal = $item[1]
    $jVal = $item[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
